$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# This workbook tracks localization handoff status for two source files:
#   4dc45583-c25c-4cc7-ba22-5017bcc1409b.md
#   f2cb6b02-2037-440e-a631-2058e246c16a.md
# A new handoff report was generated: the two rows (row 2 / row 3) swap
# places on every sheet, their "Latest Handoff Date(time)" values are
# refreshed, and 4dc45583's status moves from "Handed back: in sync with
# en-US" to "Ready for handoff".
# -----------------------------------------------------------------------

$missing = [Type]::Missing

function Set-Hyperlink {
    param($ws, [string]$cellRef, [string]$address, [string]$display)
    $rng = $ws.Range($cellRef)
    if ($rng.Hyperlinks.Count -gt 0) {
        $rng.Hyperlinks.Delete()
    }
    $ws.Hyperlinks.Add($rng, $address, $missing, $missing, $display) | Out-Null
}

# =========================================================================
# Sheet "Overview"
# =========================================================================
$ws1 = $wb.Worksheets.Item("Overview")

Set-Hyperlink $ws1 "A2" "https://github.com/OpenLocalizationTest/oltest/blob/dcfa114d60b47d912db7a8ad469578644a86f4e9/e2e/4dc45583-c25c-4cc7-ba22-5017bcc1409b.md" "f2cb6b02-2037-440e-a631-2058e246c16a.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"
$ws1.Range("D2").Value = "2016-03-19 08:45:31"

Set-Hyperlink $ws1 "A3" "https://github.com/OpenLocalizationTest/oltest/blob/dcfa114d60b47d912db7a8ad469578644a86f4e9/e2e/f2cb6b02-2037-440e-a631-2058e246c16a.md" "4dc45583-c25c-4cc7-ba22-5017bcc1409b.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-03-19 08:45:31"

# =========================================================================
# Sheet "zh-cn"
# =========================================================================
$ws2 = $wb.Worksheets.Item("zh-cn")

Set-Hyperlink $ws2 "A2" "https://github.com/OpenLocalizationTest/oltest/blob/dcfa114d60b47d912db7a8ad469578644a86f4e9/e2e/4dc45583-c25c-4cc7-ba22-5017bcc1409b.md" "f2cb6b02-2037-440e-a631-2058e246c16a.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
Set-Hyperlink $ws2 "D2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/027aca9c269530f2c41ef2367bd081df3627b3f0/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/4dc45583-c25c-4cc7-ba22-5017bcc1409b.d22e9399b183375bf9a64293cd88573e7fcd677f.zh-cn.xlf" "f2cb6b02-2037-440e-a631-2058e246c16a.8fd2d8b6642dc3fc051782fecd1ef3bf24bafd2a.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-19 08:45:22"
Set-Hyperlink $ws2 "F2" "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/4a74f69850b6a8b5970047c30bb6ce8a6bb0562f/e2e/4dc45583-c25c-4cc7-ba22-5017bcc1409b.md" "f2cb6b02-2037-440e-a631-2058e246c16a.md"
Set-Hyperlink $ws2 "G2" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8fca803eb5de61a28fd1a070e88ee632199214ac/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/4dc45583-c25c-4cc7-ba22-5017bcc1409b.d22e9399b183375bf9a64293cd88573e7fcd677f.zh-cn.xlf" "f2cb6b02-2037-440e-a631-2058e246c16a.8fd2d8b6642dc3fc051782fecd1ef3bf24bafd2a.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-03-19 08:44:19"
$ws2.Range("I2").Value = ""
$ws2.Range("J2").Value = "Include"

Set-Hyperlink $ws2 "A3" "https://github.com/OpenLocalizationTest/oltest/blob/dcfa114d60b47d912db7a8ad469578644a86f4e9/e2e/f2cb6b02-2037-440e-a631-2058e246c16a.md" "4dc45583-c25c-4cc7-ba22-5017bcc1409b.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
Set-Hyperlink $ws2 "D3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/027aca9c269530f2c41ef2367bd081df3627b3f0/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/f2cb6b02-2037-440e-a631-2058e246c16a.8fd2d8b6642dc3fc051782fecd1ef3bf24bafd2a.zh-cn.xlf" "4dc45583-c25c-4cc7-ba22-5017bcc1409b.d22e9399b183375bf9a64293cd88573e7fcd677f.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-19 08:45:22"
Set-Hyperlink $ws2 "F3" "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/4a74f69850b6a8b5970047c30bb6ce8a6bb0562f/e2e/f2cb6b02-2037-440e-a631-2058e246c16a.md" "4dc45583-c25c-4cc7-ba22-5017bcc1409b.md"
Set-Hyperlink $ws2 "G3" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8fca803eb5de61a28fd1a070e88ee632199214ac/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/f2cb6b02-2037-440e-a631-2058e246c16a.8fd2d8b6642dc3fc051782fecd1ef3bf24bafd2a.zh-cn.xlf" "4dc45583-c25c-4cc7-ba22-5017bcc1409b.d22e9399b183375bf9a64293cd88573e7fcd677f.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-03-19 08:44:19"
$ws2.Range("I3").Value = ""
$ws2.Range("J3").Value = "Include"

# =========================================================================
# Sheet "de-de"
# =========================================================================
$ws3 = $wb.Worksheets.Item("de-de")

Set-Hyperlink $ws3 "A2" "https://github.com/OpenLocalizationTest/oltest/blob/dcfa114d60b47d912db7a8ad469578644a86f4e9/e2e/4dc45583-c25c-4cc7-ba22-5017bcc1409b.md" "f2cb6b02-2037-440e-a631-2058e246c16a.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
Set-Hyperlink $ws3 "D2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2103be1f6acc92d14de0f3af83d8ed1938eb9b49/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/4dc45583-c25c-4cc7-ba22-5017bcc1409b.d22e9399b183375bf9a64293cd88573e7fcd677f.de-de.xlf" "f2cb6b02-2037-440e-a631-2058e246c16a.8fd2d8b6642dc3fc051782fecd1ef3bf24bafd2a.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-19 08:45:31"
Set-Hyperlink $ws3 "F2" "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/ce48fe823af42462157398306b658be7f775ce66/e2e/4dc45583-c25c-4cc7-ba22-5017bcc1409b.md" "f2cb6b02-2037-440e-a631-2058e246c16a.md"
Set-Hyperlink $ws3 "G2" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/bb783c65ea87b2961f5eaf2f7c387e1293ab19e2/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/4dc45583-c25c-4cc7-ba22-5017bcc1409b.d22e9399b183375bf9a64293cd88573e7fcd677f.de-de.xlf" "f2cb6b02-2037-440e-a631-2058e246c16a.8fd2d8b6642dc3fc051782fecd1ef3bf24bafd2a.de-de.xlf"
$ws3.Range("H2").Value = "2016-03-19 08:44:32"
$ws3.Range("I2").Value = ""
$ws3.Range("J2").Value = "Include"

Set-Hyperlink $ws3 "A3" "https://github.com/OpenLocalizationTest/oltest/blob/dcfa114d60b47d912db7a8ad469578644a86f4e9/e2e/f2cb6b02-2037-440e-a631-2058e246c16a.md" "4dc45583-c25c-4cc7-ba22-5017bcc1409b.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
Set-Hyperlink $ws3 "D3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2103be1f6acc92d14de0f3af83d8ed1938eb9b49/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/f2cb6b02-2037-440e-a631-2058e246c16a.8fd2d8b6642dc3fc051782fecd1ef3bf24bafd2a.de-de.xlf" "4dc45583-c25c-4cc7-ba22-5017bcc1409b.d22e9399b183375bf9a64293cd88573e7fcd677f.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-19 08:45:31"
Set-Hyperlink $ws3 "F3" "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/ce48fe823af42462157398306b658be7f775ce66/e2e/f2cb6b02-2037-440e-a631-2058e246c16a.md" "4dc45583-c25c-4cc7-ba22-5017bcc1409b.md"
Set-Hyperlink $ws3 "G3" "https://github.com/OpenLocalizationTestOrg/olhandback/blob/bb783c65ea87b2961f5eaf2f7c387e1293ab19e2/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/f2cb6b02-2037-440e-a631-2058e246c16a.8fd2d8b6642dc3fc051782fecd1ef3bf24bafd2a.de-de.xlf" "4dc45583-c25c-4cc7-ba22-5017bcc1409b.d22e9399b183375bf9a64293cd88573e7fcd677f.de-de.xlf"
$ws3.Range("H3").Value = "2016-03-19 08:44:32"
$ws3.Range("I3").Value = ""
$ws3.Range("J3").Value = "Include"

Write-Host "Handoff report regenerated."
